$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row content: "q" -> "question", "a" -> "answer"
$ws.Range("A1").Value = "question"
$ws.Range("B1").Value = "answer"

# Update the active selection to B1 (as recorded in the saved sheet view)
$ws.Range("B1").Select()
